$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 12.93898809523811
$ws.Range("N2").Value = 1.586442583591966
$ws.Range("O2").Value = 1.700608911205746
$ws.Range("I3").Value = 12.93898809523811
$ws.Range("I7").Value = -0.763888888888889
$ws.Range("N7").Value = 1.432007893438579
$ws.Range("O7").Value = 1.522400424853956
$ws.Range("I10").Value = 19.65277777777778
$ws.Range("N10").Value = 1.674945181765724
$ws.Range("O10").Value = 1.804078036500944
$ws.Range("I13").Value = 5.462962962962945
$ws.Range("N13").Value = 1.49828630419821
$ws.Range("O13").Value = 1.598520446096654
$ws.Range("I15").Value = 5.462962962962945
$ws.Range("N15").Value = 1.49828630419821
$ws.Range("O15").Value = 1.598520446096654
$ws.Range("I17").Value = 5.462962962962945
$ws.Range("N17").Value = 1.49828630419821
$ws.Range("O17").Value = 1.598520446096654
$ws.Range("I20").Value = 1.791666666666668
$ws.Range("N20").Value = 1.458486584262888
$ws.Range("O20").Value = 1.552746181345467
$ws.Range("I21").Value = 1.791666666666668
$ws.Range("I23").Value = 12.66820987654322
$ws.Range("N23").Value = 1.583068924143447
$ws.Range("O23").Value = 1.696684247214952
$ws.Range("I25").Value = 5.462962962962945
$ws.Range("N25").Value = 1.49828630419821
$ws.Range("O25").Value = 1.598520446096654
$ws.Range("I26").Value = 19.60879629629628
$ws.Range("N26").Value = 1.674333288469303
$ws.Range("O26").Value = 1.803359265239363
$ws.Range("I27").Value = 14.96875
$ws.Range("N27").Value = 1.612196950762309
$ws.Range("O27").Value = 1.730616680249932
$ws.Range("I28").Value = 14.96875
$ws.Range("I29").Value = 15.36574074074072
$ws.Range("N29").Value = 1.617332194197838
$ws.Range("O29").Value = 1.73660999151892
$ws.Range("I30").Value = 13.46442495126706
$ws.Range("N30").Value = 1.593030259848797
$ws.Range("O30").Value = 1.708276634982499
$ws.Range("I36").Value = 19.65277777777778
$ws.Range("N36").Value = 1.674945181765724
$ws.Range("O36").Value = 1.804078036500944
$ws.Range("I37").Value = 14.47727272727272
$ws.Range("N37").Value = 1.605884483070795
$ws.Range("O37").Value = 1.723253983867794
$ws.Range("I38").Value = 14.47727272727272
$ws.Range("I39").Value = 19.65277777777778
$ws.Range("N39").Value = 1.674945181765724
$ws.Range("O39").Value = 1.804078036500944
$ws.Range("I41").Value = 12.93898809523811
$ws.Range("N41").Value = 1.586442583591966
$ws.Range("O41").Value = 1.700608911205746
$ws.Range("I42").Value = 1.925925925925943
$ws.Range("N42").Value = 1.459904774678112
$ws.Range("O42").Value = 1.554373915558126
